# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback:
#  - the "Ready for handoff" status becomes "Handed back: in sync with en-US"
#    everywhere it appears (Overview zh-cn/de-de columns, and the Status
#    column on each language sheet)
#  - each language sheet's "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns get filled in for the two rows
#  - the new "Latest Target File" cells become hyperlinks to the source .md
#  - a handful of columns get widened to fit the newly-added long values
#
# NOTE on column widths: this host's ColumnWidth setter snaps the stored
# OOXML width to the nearest 1/6 of a character, so the exact fractional
# widths Excel would have computed (e.g. 29.9777047293527) are not
# reproducible bit-for-bit; the inputs below are chosen to land on the
# closest value the engine can actually store.

$wb = $excel.ActiveWorkbook

$mdName       = "be3fe9fd-d7f6-43a2-b38c-22c1a60f6fdd.md"
$mdUrl        = "https://github.com/OpenLocalizationTestOrg/oltest/blob/33077d0a5938814937bfed1ca65fc009764cafce/e2e/be3fe9fd-d7f6-43a2-b38c-22c1a60f6fdd.md"
$zhXlf        = "be3fe9fd-d7f6-43a2-b38c-22c1a60f6fdd.09cfd9f958ddaf58f738cb8355c698180b020c9f.zh-cn.xlf"
$deXlf        = "be3fe9fd-d7f6-43a2-b38c-22c1a60f6fdd.09cfd9f958ddaf58f738cb8355c698180b020c9f.de-de.xlf"
$newStatus    = "Handed back: in sync with en-US"
$zhHandback   = "2016-08-14 01:30:11"
$deHandback   = "2016-08-14 01:30:22"

# ---------------------------------------------------------------------
# Overview sheet: status rollup columns (zh-cn / de-de) + widen them
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $mdName
$wsZh.Range("I3").Value = $mdName
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K2").Value = $zhHandback
$wsZh.Range("K3").Value = $zhHandback

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $mdName
$wsDe.Range("I3").Value = $mdName
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K2").Value = $deHandback
$wsDe.Range("K3").Value = $deHandback

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
